$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 427, shifting existing rows 427:465 down to 428:466
$ws.Rows.Item(427).Insert()

# Populate the newly inserted row 427 with the new weekly record
$ws.Range("A427").Value = 9
$ws.Range("B427").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C427").Value = "Metropolitana"
$ws.Range("D427").Value = 45132
$ws.Range("E427").Value = 13
$ws.Range("F427").Value = 300000001
$ws.Range("G427").Value = "Rabanito"
$ws.Range("H427").Value = "Sin especificar"
$ws.Range("I427").Value = "Primera"
$ws.Range("J427").Value = 7000
$ws.Range("K427").Value = 3000
$ws.Range("L427").Value = 4000
$ws.Range("M427").Value = 3500
$ws.Range("N427").Value = "`$/cien unidades (volumen en unidades)"
$ws.Range("O427").Value = "Región Metropolitana"
$ws.Range("P427").Value = 35
$ws.Range("Q427").Value = 100
$ws.Range("R427").Value = "Hortaliza"
